$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Status text update: "Ready for handoff" -> "In Translation"
#    This shared string is used by:
#      Overview!E2, Overview!F2, Overview!E3, Overview!F3  (zh-cn / de-de cols)
#      zh-cn!C2,   zh-cn!C3                                (Status col)
#      de-de!C2,   de-de!C3                                (Status col)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# ---------------------------------------------------------------------------
# 2) Narrow the "Status" columns now that the text is shorter:
#      Overview columns E and F (zh-cn / de-de) : ~17.22 -> ~13.41 chars
#      zh-cn column C (Status)                  : ~17.22 -> ~13.41 chars
#      de-de column C (Status)                  : ~17.22 -> ~13.41 chars
#    ColumnWidth is expressed in "characters"; 12.5 lands squarely in the
#    pixel bucket closest to the target width.
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
